# Apply the edits described by the diff:
#  - Sheet "maxRel": A1 value 20000 -> 25000; selection set to F11
#  - Sheet "minRel": A1 value 8000 -> 5000; selection set to G12

$wb = $excel.ActiveWorkbook

# --- maxRel sheet ---
$wsMaxRel = $wb.Worksheets.Item("maxRel")
$wsMaxRel.Range("A1").Value = 25000
$wsMaxRel.Activate()
$wsMaxRel.Range("F11").Select()

# --- minRel sheet ---
$wsMinRel = $wb.Worksheets.Item("minRel")
$wsMinRel.Range("A1").Value = 5000
$wsMinRel.Activate()
$wsMinRel.Range("G12").Select()

# Restore maxRel as the active/selected tab (tabSelected="1" in diff)
$wsMaxRel.Activate()
